$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-looking string into A69 as literal text (matching the
# other rows' inlineStr "YYYY/MM/DD" cells) instead of letting Excel
# auto-convert it to a date serial number, then drop the leftover
# text-number-format so the cell keeps the sheet's default (no explicit
# style), same as every other data row.
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "2025/10/06"
$ws.Range("A69").ClearFormats()

$ws.Range("B69").Value = "月"
$ws.Range("C69").Value = 16
$ws.Range("D69").Value = 61
